# The deck's slides 7 and 8 (SlideID 272 = "Результат" recap slide, and
# SlideID 278 = "Рисунок 5" / PSB regional network slide) were swapped in
# display order. Reproduce that by moving slide 7 to position 8 (which
# pushes the former slide 8 up to position 7).
$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(7)
$slide.MoveTo(8)
